$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 69

# Update row 3 values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 30

# Delete the former row 4 entirely, shrinking the used range to A1:B3
$ws.Rows("4:4").Delete()
